$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '26.544.55'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +0.13%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.814.31'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +0.17%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.004'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.19%  '
$ws.Range('E5').Value = '  -0.20%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '305.83'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.69%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4546'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.15%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3594'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -1.91%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '46.30'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +2.79%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.07108'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.39%  '
$ws.Range('E11').Value = '  +1.46%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.07722'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.50%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '19.31'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -0.30%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '1.840.51'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +1.71%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '5.258'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.56%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '6.301'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -1.00%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '85.80'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.68%  '
$ws.Range('E18').Value = '  -0.22%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.000008541'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.69%  '
$ws.Range('E20').Value = '  -0.23%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '26.587.92'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.03%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '14.15'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.68%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.954'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.97%  '
$ws.Range('E24').Value = '  +0.33%  '
$ws.Range('E25').Value = '  -3.18%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '152.01'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.26%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '17.77'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.96%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.019'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -1.96%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '112.21'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.60%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.816'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.80%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.08714'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.39%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.135'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +2.48%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.7449'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +1.44%  '
$ws.Range('B34').Value = 'Filecoin'
$ws.Range('C34').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.430'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -1.74%  '
$ws.Range('B35').Value = 'RenderToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.719'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +1.20%  '
$ws.Range('E36').Value = '  -0.90%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.072'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.89%  '
$ws.Range('E38').Value = '  -0.70%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.912'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.27%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.05081'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.58%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.5091'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +1.74%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '6.787'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -2.90%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.1505'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -3.49%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '8.021'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -1.88%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.4686'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +1.69%  '
$ws.Range('E46').Value = '  -0.25%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '9.931'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.10%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '98.85'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -2.06%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.561'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -1.87%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.05995'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.01%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '63.80'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.94%  '
